# Applies the cryptos list update described in the commit: price/volume refresh
# for Sat Feb 17 18:32:13 UTC 2024, including a Chainlink/Dogecoin row swap (rows 12-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.249.17"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "'2.771.06"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'353.42"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "'107.57"
$ws.Range("E6").Value = "  -1.29%  "
$ws.Range("D7").Value = "'0.549"
$ws.Range("E7").Value = "  -2.18%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "'39.42"
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("E11").Value = "  +3.36%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "'19.93"
$ws.Range("E12").Value = "  +2.93%  "
$ws.Range("B13").Value = "Dogecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D13").Value = "'0.0830"
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "'3.203.31"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "'2.758.15"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "'0.925"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "'51.144.52"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("E19").Value = "  +3.49%  "
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "'13.09"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").Value = "'69.54"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "'265.41"
$ws.Range("E24").Value = "  -3.27%  "
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'25.94"
$ws.Range("E28").Value = "  +12.51%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").Value = "'35.52"
$ws.Range("E31").Value = "  +5.28%  "
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").Value = "'6.06"
$ws.Range("E33").Value = "  +6.43%  "
$ws.Range("E34").Value = "  -5.08%  "
$ws.Range("D35").Value = "'5.53"
$ws.Range("E35").Value = "  +4.32%  "
$ws.Range("D36").Value = "'0.0826"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'18.12"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").Value = "'1.95"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "'22.07"
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("D45").Value = "'2.19"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "'2.095.08"
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").Value = "'0.906"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("D50").Value = "'5.41"
$ws.Range("E50").Value = "  -4.97%  "
$ws.Range("E51").Value = "  +6.78%  "
